$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2435.7144
$ws.Range("I40").Value = 1762
$ws.Range("J40").Value = 3334
$ws.Range("K40").Value = 1762
$ws.Range("L40").Value = 3334
$ws.Range("M40").Value = -1587
$ws.Range("N40").Value = -3684
$ws.Range("H76").Value = 4114
$ws.Range("J76").Value = 4114
$ws.Range("L76").Value = 4114
$ws.Range("N76").Value = -4744
$ws.Range("H79").Value = 4114
$ws.Range("J79").Value = 4114
$ws.Range("L79").Value = 4114
$ws.Range("N79").Value = -6298
$ws.Range("H121").Value = 3933.3333
$ws.Range("J121").Value = 3933.3333
$ws.Range("L121").Value = 11799.9999
$ws.Range("N121").Value = -15293.9999
$ws.Range("H133").Value = 89779
$ws.Range("J133").Value = 89779
$ws.Range("L133").Value = 89779
$ws.Range("N133").Value = -99899
$ws.Range("H137").Value = 70521.89
$ws.Range("I137").Value = 63518.2
$ws.Range("K137").Value = 190554.6
$ws.Range("M137").Value = -188004.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21175.27
$ws.Range("I32").Value = 21742.3
$ws.Range("J32").Value = 6999.5
$ws.Range("K32").Value = 21742.3
$ws.Range("L32").Value = 6999.5
$ws.Range("M32").Value = -21455.3
$ws.Range("N32").Value = -7573.5
$ws.Range("H61").Value = 4936.625
$ws.Range("I61").Value = 1319.2084
$ws.Range("K61").Value = 1319.2084
$ws.Range("M61").Value = -1107.2084
$ws.Range("H63").Value = 3999.8
$ws.Range("I63").Value = 1749.5
$ws.Range("K63").Value = 1749.5
$ws.Range("M63").Value = -1063.5
$ws.Range("H66").Value = 3999.8
$ws.Range("I66").Value = 1749.5
$ws.Range("K66").Value = 8747.5
$ws.Range("M66").Value = -5315.5
$ws.Range("H74").Value = 437543.34
$ws.Range("I74").Value = 750836.75
$ws.Range("K74").Value = 750836.75
$ws.Range("M74").Value = -749962.75
$ws.Range("H77").Value = 437543.34
$ws.Range("I77").Value = 750836.75
$ws.Range("K77").Value = 3754183.75
$ws.Range("M77").Value = -3749815.75
$ws.Range("H88").Value = 14245.25
$ws.Range("I88").Value = 2490.5
$ws.Range("J88").Value = 26000
$ws.Range("K88").Value = 2490.5
$ws.Range("L88").Value = 26000
$ws.Range("M88").Value = -2084.5
$ws.Range("N88").Value = -26812
$ws.Range("H91").Value = 14245.25
$ws.Range("I91").Value = 2490.5
$ws.Range("J91").Value = 26000
$ws.Range("K91").Value = 2490.5
$ws.Range("L91").Value = 26000
$ws.Range("M91").Value = -1086.5
$ws.Range("N91").Value = -28808
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H102").Value = 2054.6
$ws.Range("I102").Value = 2054.6
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2054.6
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -432.5999999999999
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 2492.6155
$ws.Range("I132").Value = 1489.4445
$ws.Range("K132").Value = 4468.333500000001
$ws.Range("M132").Value = -1938.333500000001
$ws.Range("H136").Value = 4936.625
$ws.Range("I136").Value = 1319.2084
$ws.Range("K136").Value = 3957.6252
$ws.Range("M136").Value = -1407.6252

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2581.6667
$ws.Range("I86").Value = 2498
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2498
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -1375
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 2581.6667
$ws.Range("I89").Value = 2498
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 12490
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -6874
$ws.Range("N89").Value = -26232
$ws.Range("H122").Value = 120000
$ws.Range("J122").Value = 120000
$ws.Range("L122").Value = 120000
$ws.Range("N122").Value = -129800
$ws.Range("H134").Value = 2826.2307
$ws.Range("I134").Value = 2669.913
$ws.Range("K134").Value = 8009.739
$ws.Range("M134").Value = -5474.739

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5666.6665
$ws.Range("I6").Value = 5666.6665
$ws.Range("K6").Value = 5666.6665
$ws.Range("M6").Value = -5553.6665
$ws.Range("H22").Value = 1375
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1375
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1375
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2075
$ws.Range("H31").Value = 20009280
$ws.Range("I31").Value = 33344300
$ws.Range("K31").Value = 33344300
$ws.Range("M31").Value = -33344005
$ws.Range("H34").Value = 20009280
$ws.Range("I34").Value = 33344300
$ws.Range("K34").Value = 33344300
$ws.Range("M34").Value = -33344098
$ws.Range("H132").Value = 144144.14
$ws.Range("I132").Value = 144144.14
$ws.Range("K132").Value = 432432.42
$ws.Range("M132").Value = -429902.42
$ws.Range("H134").Value = 2193.0715
$ws.Range("I134").Value = 1900.3846
$ws.Range("K134").Value = 5701.1538
$ws.Range("M134").Value = -3166.1538

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 152.90909
$ws.Range("I14").Value = 152.90909
$ws.Range("K14").Value = 458.72727
$ws.Range("M14").Value = -285.72727
$ws.Range("H23").Value = 429.1
$ws.Range("J23").Value = 674.25
$ws.Range("L23").Value = 2022.75
$ws.Range("N23").Value = -2492.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11996.667
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H70").Value = 5076
$ws.Range("I70").Value = 4406.3335
$ws.Range("J70").Value = 6281.4
$ws.Range("K70").Value = 4406.3335
$ws.Range("L70").Value = 6281.4
$ws.Range("M70").Value = -4136.3335
$ws.Range("N70").Value = -6821.4
$ws.Range("H73").Value = 5076
$ws.Range("I73").Value = 4406.3335
$ws.Range("J73").Value = 6281.4
$ws.Range("K73").Value = 4406.3335
$ws.Range("L73").Value = 6281.4
$ws.Range("M73").Value = -3470.3335
$ws.Range("N73").Value = -8153.4
$ws.Range("H93").Value = 100000
$ws.Range("J93").Value = 100000
$ws.Range("L93").Value = 100000
$ws.Range("N93").Value = -103744
$ws.Range("H102").Value = 3715.7144
$ws.Range("I102").Value = 2753
$ws.Range("K102").Value = 2753
$ws.Range("M102").Value = -1131
$ws.Range("H117").Value = 51310
$ws.Range("J117").Value = 51310
$ws.Range("L117").Value = 51310
$ws.Range("N117").Value = -58194
$ws.Range("H132").Value = 3473.6
$ws.Range("I132").Value = 2935.7778
$ws.Range("K132").Value = 8807.3334
$ws.Range("M132").Value = -6277.3334
$ws.Range("H136").Value = 98946.39999999999
$ws.Range("J136").Value = 98946.39999999999
$ws.Range("L136").Value = 296839.2
$ws.Range("N136").Value = -301939.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4899.5
$ws.Range("I40").Value = 4899.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4899.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4763.5
$ws.Range("N40").ClearContents()
$ws.Range("H93").Value = 1300
$ws.Range("I93").Value = 1489
$ws.Range("J93").Value = 449.5
$ws.Range("K93").Value = 1489
$ws.Range("L93").Value = 449.5
$ws.Range("M93").Value = -241
$ws.Range("N93").Value = -2945.5
$ws.Range("H95").Value = 24499.5
$ws.Range("J95").Value = 24499.5
$ws.Range("L95").Value = 24499.5
$ws.Range("N95").Value = -29991.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 5134.4707
$ws.Range("I132").Value = 4935.1816
$ws.Range("K132").Value = 14805.5448
$ws.Range("M132").Value = -12275.5448

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16347.848
$ws.Range("I136").Value = 18881.846
$ws.Range("J136").Value = 2229.8572
$ws.Range("K136").Value = 56645.538
$ws.Range("L136").Value = 6689.571599999999
$ws.Range("M136").Value = -54095.538
$ws.Range("N136").Value = -11789.5716
